$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = -10.4029
$ws.Range("C18").Value = -12.84499999999999
$ws.Range("C20").Value = -11.999
$ws.Range("C27").Value = -12.7436
$ws.Range("C69").Value = -11.45259999999999
$ws.Range("C76").Value = -12.0788
$ws.Range("C82").Value = -11.84899999999999
